$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Relations sheet: add three new rows describing the "undeclared" relation
# (an error-message example): undeclared, v = D, w = E
# ---------------------------------------------------------------------------
$wsRelations = $wb.Worksheets.Item("Relations")
$wsRelations.Range("A6").Value2 = "undeclared"

$wsRelations.Range("A7").Value2 = "v"
$wsRelations.Range("B7").Value2 = "D"
$wsRelations.Range("C7").Value2 = "D"

$wsRelations.Range("A8").Value2 = "w"
$wsRelations.Range("B8").Value2 = "E"
$wsRelations.Range("C8").Value2 = "E"

# ---------------------------------------------------------------------------
# Rules sheet: add two new rule rows and widen the first two columns so the
# new, longer rule text fits
# ---------------------------------------------------------------------------
$wsRules = $wb.Worksheets.Item("Rules")

$wsRules.Range("A4").Value2 = "v = w"
$wsRules.Range("B4").Value2 = "v"
$wsRules.Range("C4").Value2 = "w"

$wsRules.Range("D5").Value2 = "undeclared = s;t"
$wsRules.Range("E5").Value2 = "undeclared"
$wsRules.Range("F5").Value2 = "s;t"

$wsRules.Columns.Item(1).ColumnWidth = 15.6
$wsRules.Columns.Item(2).ColumnWidth = 11.0

# ---------------------------------------------------------------------------
# Compositions sheet: move the "r;s = r;s" composition example from E4:G4
# up/left to A4:C4
# ---------------------------------------------------------------------------
$wsCompositions = $wb.Worksheets.Item("Compositions")
$wsCompositions.Range("A4").Value2 = $wsCompositions.Range("E4").Value2
$wsCompositions.Range("B4").Value2 = $wsCompositions.Range("F4").Value2
$wsCompositions.Range("C4").Value2 = $wsCompositions.Range("G4").Value2
$wsCompositions.Range("E4:G4").Clear()

$excel.ActiveWindow.Zoom = 160
$null = $wsCompositions.Range("A4:C4").Select()

# ---------------------------------------------------------------------------
# View state: Relations becomes the active / selected sheet, with a new
# selection further down reflecting the added rows; Rules keeps its own
# updated selection
# ---------------------------------------------------------------------------
$null = $wsRules.Range("C5").Select()

$wsRelations.Activate()
$null = $wsRelations.Range("C10").Select()
